$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

Set-TextCell 2 4 '42.850.76'
Set-TextCell 2 5 '  +0.38%  '
Set-TextCell 3 4 '2.533.65'
Set-TextCell 3 5 '  +0.60%  '
Set-TextCell 4 5 '  -0.07%  '
Set-TextCell 5 4 '315.42'
Set-TextCell 5 5 '  -0.64%  '
Set-TextCell 6 4 '96.15'
Set-TextCell 6 5 '  +0.46%  '
Set-TextCell 7 4 '0.576'
Set-TextCell 7 5 '  -1.50%  '
Set-TextCell 9 5 '  -0.93%  '
Set-TextCell 10 4 '36.20'
Set-TextCell 10 5 '  +0.30%  '
Set-TextCell 11 5 '  -0.50%  '
Set-TextCell 12 4 '7.57'
Set-TextCell 12 5 '  -2.31%  '
Set-TextCell 13 5 '  -2.75%  '
Set-TextCell 14 4 '2.921.73'
Set-TextCell 14 5 '  +0.51%  '
Set-TextCell 15 4 '2.553.98'
Set-TextCell 15 5 '  +1.32%  '
Set-TextCell 16 5 '  -1.77%  '
Set-TextCell 17 5 '  -0.96%  '
Set-TextCell 18 4 '42.891.62'
Set-TextCell 18 5 '  +0.51%  '
Set-TextCell 19 4 '13.11'
Set-TextCell 19 5 '  +1.52%  '
Set-TextCell 20 4 '6.84'
Set-TextCell 20 5 '  +3.71%  '
Set-TextCell 21 4 '0.0₃0965'
Set-TextCell 21 5 '  -1.11%  '
Set-TextCell 22 4 '70.03'
Set-TextCell 22 5 '  -2.13%  '
Set-TextCell 23 4 '253.47'
Set-TextCell 23 5 '  +0.04%  '
Set-TextCell 24 4 '2.95'
Set-TextCell 24 5 '  -1.37%  '
Set-TextCell 25 5 '  +1.70%  '
Set-TextCell 26 4 '26.66'
Set-TextCell 26 5 '  -1.65%  '
Set-TextCell 27 5 '  -0.03%  '
Set-TextCell 28 5 '  +2.47%  '
Set-TextCell 29 4 '40.56'
Set-TextCell 29 5 '  +6.70%  '
Set-TextCell 30 4 '10.44'
Set-TextCell 30 5 '  +2.70%  '
Set-TextCell 31 4 '5.94'
Set-TextCell 31 5 '  +0.29%  '
Set-TextCell 32 4 '158.01'
Set-TextCell 32 5 '  +1.46%  '
Set-TextCell 33 4 '2.17'
Set-TextCell 33 5 '  +3.75%  '
Set-TextCell 34 2 'Celestia'
Set-TextCell 34 3 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell 34 4 '19.16'
Set-TextCell 34 5 '  -1.43%  '
Set-TextCell 35 2 'LidoDAOToken'
Set-TextCell 35 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 35 4 '3.34'
Set-TextCell 35 5 '  -0.49%  '
Set-TextCell 36 5 '  +2.04%  '
Set-TextCell 37 4 '0.0783'
Set-TextCell 37 5 '  -0.82%  '
Set-TextCell 38 5 '  -1.43%  '
Set-TextCell 39 5 '  -1.38%  '
Set-TextCell 40 5 '  -4.65%  '
Set-TextCell 41 5 '  +13.64%  '
Set-TextCell 42 4 '3.85'
Set-TextCell 42 5 '  -0.44%  '
Set-TextCell 43 5 '  +0.11%  '
Set-TextCell 44 5 '  +0.22%  '
Set-TextCell 45 5 '  -2.39%  '
Set-TextCell 46 4 '2.032.41'
Set-TextCell 46 5 '  +0.04%  '
Set-TextCell 47 4 '9.21'
Set-TextCell 47 5 '  +2.51%  '
Set-TextCell 48 4 '85.15'
Set-TextCell 48 5 '  +0.62%  '
Set-TextCell 49 4 '106.35'
Set-TextCell 49 5 '  +4.58%  '
Set-TextCell 50 4 '74.84'
Set-TextCell 50 5 '  -0.11%  '
Set-TextCell 51 4 '2.775.74'
Set-TextCell 51 5 '  +0.47%  '
